# Fix not appending of excel:
#  - "CSVworks" row: name column CSV_FAIL -> CSV_PASS
#  - "ExcelButNoSets" row: name column EXCEL_FAIL -> EXCEL_PASS
#  - "ExcelButNoSets" row: partID column DoesntExistInSets -> ExistInSets

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("files")

$ws.Range("C2").Value = "CSV_PASS"
$ws.Range("C3").Value = "EXCEL_PASS"
$ws.Range("E3").Value = "ExistInSets"
